# Add a new "Coin" item row to the Entities master-item sheet and
# update the remembered selection, matching the upstream commit's
# "MasterQuest and UI quest show" data addition.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# New row (id=9, key="Coin", value=10304) appended right after the
# existing last data row (row 9).
$newRow = 10
$ws.Range("A$newRow").Value = 9
$ws.Range("B$newRow").Value = "Coin"
$ws.Range("C$newRow").Value = 10304

# Carry over the same row formatting ("Neutral" cell style, s="3") used
# by the previous data row so the new row renders consistently.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Reflect the author's last active selection recorded in the sheet.
$ws.Range("D12").Select()
